$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Germany)
$ws.Range("C2").Value = 0.06653805445417677

# Row 3 (Spain)
$ws.Range("C3").Value = 0.2179061588686753
$ws.Range("D3").Value = -0.1513681044144985
$ws.Range("E3").Value = -0.0767665620947065
$ws.Range("F3").Value = -0.07460154231979201
$ws.Range("G3").Value = 0.5071515058713623
$ws.Range("H3").Value = 0.4928484941286379

# Row 4 (Italy)
$ws.Range("C4").Value = 0.2250579911111283
$ws.Range("D4").Value = -0.1585199366569516
$ws.Range("E4").Value = -0.08886230598376316
$ws.Range("F4").Value = -0.06965763067318842
$ws.Range("G4").Value = 0.5605749526387177
$ws.Range("H4").Value = 0.4394250473612823
